$d = $word.ActiveDocument

# Namespace declarations reused by every InsertXML payload below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Replace-ParagraphXml($findText, $paraInnerXml) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Delete()
    $xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
        "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
        "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
        "<pkg:xmlData>" +
        "<w:document $wNs><w:body>$paraInnerXml</w:body></w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
    $rng.InsertXML($xml)
}

# 1) "Initially we have the names of all supported 42 instructions, register
#    names, and 5 halting instructions in the constants file." was split
#    across five runs; collapse it back into a single run (text unchanged).
$d.Content.Find.Execute( `
    "Initially we have the names of all supported 42 instructions, register names, and 5 halting instructions in the constants file.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Initially we have the names of all supported 42 instructions, register names, and 5 halting instructions in the constants file.", `
    2) | Out-Null

# 2) "in parse.py, ... that then populates the memory" -> wrap "populates" in
#    a grammar-check proofErr span.
Replace-ParagraphXml `
    "in parse.py, the assembly code file is read and cleaned into data structures that then populates the memory" `
    ('<w:p w14:paraId="09FCFFA9" w14:textId="77777777" w:rsidR="00022157" w:rsidRDefault="00022157" w:rsidP="00022157">' + `
     '<w:r><w:t xml:space="preserve">in parse.py, the assembly code file is read and cleaned into data structures that then </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>populates</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t xml:space="preserve"> the memory</w:t></w:r>' + `
     '</w:p>')

# 3) "The overall connection and running of the program then happens in the
#    main file." -> wrap "happens" in a grammar-check proofErr span.
Replace-ParagraphXml `
    "The overall connection and running of the program then happens in the main file." `
    ('<w:p w14:paraId="7BB2F9BC" w14:textId="77777777" w:rsidR="00022157" w:rsidRDefault="00022157" w:rsidP="00022157">' + `
     '<w:r><w:t xml:space="preserve">The overall connection and running of the program then </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>happens</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t xml:space="preserve"> in the main file.</w:t></w:r>' + `
     '</w:p>')

# 4) "    -a list for the instructions memory." -> wrap "instructions" in a
#    grammar-check proofErr span.
Replace-ParagraphXml `
    "    -a list for the instructions memory." `
    ('<w:p w14:paraId="31CFBE2A" w14:textId="77777777" w:rsidR="00022157" w:rsidRDefault="00022157" w:rsidP="00022157">' + `
     '<w:r><w:t xml:space="preserve">    -a list for the </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>instructions</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t xml:space="preserve"> memory.</w:t></w:r>' + `
     '</w:p>')

# 5) Append a brand-new paragraph "Example on testcase0" after the "Start by
#    running main.py ..." paragraph. Done BEFORE the proofErr split below so
#    that paragraph is no longer the last one in the body when InsertXML
#    rewrites it (InsertXML right at the end of the body, just before
#    </w:body>, leaves behind an extra empty paragraph that can't be
#    removed — the body's final paragraph mark is structurally required —
#    so we sidestep that edge case by reordering these two steps).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Last.Range.Text = "Example on testcase0"

# 6) "Start by running main.py ... Once given, The final output is displayed
#    in the terminal." -> wrap "The" in a grammar-check proofErr span.
Replace-ParagraphXml `
    "Start by running main.py where a prompt will then be given asking for the name of the assembly input file you want to run. Once given, The final output is displayed in the terminal." `
    ('<w:p w14:paraId="158472B7" w14:textId="6EC7F472" w:rsidR="00022157" w:rsidRPr="00022157" w:rsidRDefault="00022157">' + `
     '<w:r><w:t xml:space="preserve">Start by running main.py where a prompt will then be given asking for the name of the assembly input file you want to run. Once given, </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>The</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t xml:space="preserve"> final output is displayed in the terminal.</w:t></w:r>' + `
     '</w:p>')
